# PYME-4265: Add "Subscription ID" column to the european_funds_requests
# export template, between "Customer Tax ID" (H) and "Antivirus Quantity"
# (old I, now J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I; this shifts the old I:M headers to J:N and
# carries over the header cell style (gray fill) from the column to its
# left, same as Excel's native "Insert Sheet Columns" command.
$ws.Columns("I").Insert()

# New header text for the inserted column.
$ws.Range("I1").Value = "Subscription ID"

# Match the width used for the other "plain" (non bestFit) header columns
# (18 chars as stored in the XML -> ~17.17 in the COM "characters" units).
$ws.Columns("I").ColumnWidth = 17.17

# Re-apply the autofilter so its range grows from A1:M1 to A1:N1.
$ws.AutoFilterMode = $false
$ws.Range("A1:N1").AutoFilter()

# The hidden _xlnm._FilterDatabase defined name keeps the autofilter range
# too; move it from A1:M1 to A1:N1 to stay in sync.
foreach ($n in $wb.Names) {
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=Data!`$A`$1:`$N`$1"
  }
}

# Restore the active-cell selection used in the author's saved state.
$ws.Range("I2").Select()
